$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

# Update the bounding-box coordinates (Res_Bounding_Box / SK_Bounding_Box columns)
# to reflect the re-projected (WGS 84, EPSG 4269) shapefile coordinates.
$ws.Range("A2").Value = -70.9798
$ws.Range("B2").Value = -70.9658

$ws.Range("A3").Value = 42.9614
$ws.Range("B3").Value = 42.9716

$ws.Range("A4").Value = -70.7954
$ws.Range("B4").Value = -70.8095

$ws.Range("A5").Value = 43.199
$ws.Range("B5").Value = 43.1888
